$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17 currently holds the "old" weekly record (date 2021-09-27 / serial 44466).
# The edit adds a newer weekly record (date 2021-10-22 / serial 44491, with
# updated price/origin figures) as the new row 17, and pushes the previous
# record down to a brand-new row 18, unchanged.

# 1) Build the new row 18 as an exact copy of the current row 17 values/format
#    (reads use Value2, since Value getter is unreliable for COM cell reads here)
$ws.Cells.Item(18, 1).Value = $ws.Cells.Item(17, 1).Value2    # A Mercado ID
$ws.Cells.Item(18, 2).Value = $ws.Cells.Item(17, 2).Value2    # B Mercado
$ws.Cells.Item(18, 3).Value = $ws.Cells.Item(17, 3).Value2    # C Region
$ws.Cells.Item(18, 4).Value = $ws.Cells.Item(17, 4).Value2    # D Fecha
$ws.Cells.Item(18, 4).NumberFormat = $ws.Cells.Item(17, 4).NumberFormat
$ws.Cells.Item(18, 5).Value = $ws.Cells.Item(17, 5).Value2    # E Codreg
$ws.Cells.Item(18, 6).Value = $ws.Cells.Item(17, 6).Value2    # F Categoria ID
$ws.Cells.Item(18, 7).Value = $ws.Cells.Item(17, 7).Value2    # G Categoria
$ws.Cells.Item(18, 8).Value = $ws.Cells.Item(17, 8).Value2    # H Variedad
$ws.Cells.Item(18, 9).Value = $ws.Cells.Item(17, 9).Value2    # I Calidad
$ws.Cells.Item(18, 10).Value = $ws.Cells.Item(17, 10).Value2  # J Volumen
$ws.Cells.Item(18, 11).Value = $ws.Cells.Item(17, 11).Value2  # K Precio minimo
$ws.Cells.Item(18, 12).Value = $ws.Cells.Item(17, 12).Value2  # L Precio maximo
$ws.Cells.Item(18, 13).Value = $ws.Cells.Item(17, 13).Value2  # M Precio promedio ponderado
$ws.Cells.Item(18, 14).Value = $ws.Cells.Item(17, 14).Value2  # N Unidad de comercializacion
$ws.Cells.Item(18, 15).Value = $ws.Cells.Item(17, 15).Value2  # O Origen
$ws.Cells.Item(18, 16).Value = $ws.Cells.Item(17, 16).Value2  # P Precio $/Kg
$ws.Cells.Item(18, 17).Value = $ws.Cells.Item(17, 17).Value2  # Q Kg o Unidades
$ws.Cells.Item(18, 18).Value = $ws.Cells.Item(17, 18).Value2  # R Clasificacion

# 2) Overwrite row 17 in place with the newer weekly values
$ws.Cells.Item(17, 4).Value = 44491    # D17 Fecha
$ws.Cells.Item(17, 11).Value = 8000    # K17 Precio minimo
$ws.Cells.Item(17, 12).Value = 9000    # L17 Precio maximo
$ws.Cells.Item(17, 13).Value = 8500    # M17 Precio promedio ponderado
$ws.Cells.Item(17, 15).Value = "Región del Maule"  # O17 Origen
$ws.Cells.Item(17, 16).Value = 340     # P17 Precio $/Kg

$wb.Save()
